$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F4").Value = 5875
$ws.Range("F5").Value = 71
$ws.Range("F9").Value = 1566
$ws.Range("F11").Value = 29
$ws.Range("F12").Value = 672
$ws.Range("F13").Value = 1578
$ws.Range("F14").Value = 1578
$ws.Range("F15").Value = 1536
$ws.Range("F16").Value = 548
$ws.Range("F17").Value = 141
$ws.Range("F18").Value = 610
$ws.Range("F19").Value = 4395
$ws.Range("F23").Value = 812
$ws.Range("F24").Value = 3
$ws.Range("F26").Value = 2298
$ws.Range("F31").Value = 1224
$ws.Range("F32").Value = 785
$ws.Range("F34").Value = 1192
$ws.Range("F35").Value = 1184
$ws.Range("F36").Value = 80

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F12").Value = 104
$ws.Range("F13").Value = 215
$ws.Range("F15").Value = 34
$ws.Range("F18").Value = 124
$ws.Range("F19").Value = 297
$ws.Range("F21").Value = 489

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 249
$ws.Range("F3").Value = 632
$ws.Range("F4").Value = 171
$ws.Range("F5").Value = 258

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 249
$ws.Range("F6").Value = 632
$ws.Range("F7").Value = 171
$ws.Range("F8").Value = 5875
$ws.Range("F10").Value = 71
$ws.Range("F20").Value = 1566
$ws.Range("F23").Value = 29
$ws.Range("F24").Value = 1578
$ws.Range("F25").Value = 104
$ws.Range("F26").Value = 1536
$ws.Range("F27").Value = 548
$ws.Range("F28").Value = 141
$ws.Range("F29").Value = 610
$ws.Range("F30").Value = 4395
$ws.Range("F33").Value = 812
$ws.Range("F35").Value = 2298
$ws.Range("F40").Value = 1224
$ws.Range("F42").Value = 124
$ws.Range("F43").Value = 297
$ws.Range("F45").Value = 489
$ws.Range("F46").Value = 785
$ws.Range("F50").Value = 80

